# Updated symbol list on Mon Dec 19 09:44:34 UTC 2022 with GitHub Actions
#
# This script refreshes the "Price" (column D) figures for a set of coins,
# tweaks a couple of "Volume(1h)" (column E) labels, and swaps the two rows
# that list BKEXToken / CEJI (their rank prefix in column A / date columns
# F,G stay put - only B,C,D,E move).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Column D holds numeric-looking strings ("247.58", "0.03170", ...) that
    # must stay literal text (trailing zeros / exact digit count matter).
    # Pre-marking the cell as Text ("@") keeps Excel/COM from silently
    # re-interpreting the assignment as a Double and normalising it.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# ---- simple price refreshes (column D) ----
Set-TextValue "D2"  "247.58"
Set-TextValue "D3"  "21.76"
Set-TextValue "D4"  "5.464"
Set-TextValue "D5"  "0.05692"
Set-TextValue "D6"  "3.381"
Set-TextValue "D7"  "0.8048"
Set-TextValue "D8"  "1.038"
Set-TextValue "D9"  "0.1463"
Set-TextValue "D10" "0.07333"
Set-TextValue "D11" "0.03170"
Set-TextValue "D12" "0.02942"
Set-TextValue "D14" "0.001647"
Set-TextValue "D15" "3.380"
Set-TextValue "D16" "0.04726"

# row 17 also gets a Volume(1h) label tweak
Set-TextValue "D17" "0.0005859"
$ws.Range("E17").Value = "16OneONEWorstin24h"

Set-TextValue "D18" "0.006314"
Set-TextValue "D19" "0.005048"
Set-TextValue "D20" "0.001045"
Set-TextValue "D22" "0.0003200"
Set-TextValue "D24" "6.430"
Set-TextValue "D25" "2.115"
Set-TextValue "D26" "0.3291"
Set-TextValue "D40" "0.04110"
Set-TextValue "D41" "0.006961"

# ---- row 42 / row 43 swap: BKEXToken <-> CEJI ----
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003501"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1043"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue "D44" "0.008142"
Set-TextValue "D45" "0.00005813"

# row 47 also loses the "Worstin24h" suffix on its Volume(1h) label
Set-TextValue "D47" "0.0005499"
$ws.Range("E47").Value = "46ACDXExchangeACXT"

Set-TextValue "D48" "0.6824"
Set-TextValue "D49" "0.009611"
Set-TextValue "D50" "0.00002101"
